$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 20 with the new day's data
$ws.Range("A20").Value = 45964
$ws.Range("B20").Value = 5596
$ws.Range("C20").Value = 4324
$ws.Range("D20").Value = 3950
$ws.Range("E20").Value = 287
$ws.Range("F20").Value = 54
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = 6
$ws.Range("I20").Value = 0

# Update the selection to reflect the new last-entered row
$ws.Range("A20:I20").Select()
